$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 370.810666666667
$ws.Range("H3").Value = 3763.789333333334
$ws.Range("C4").Value = 143.7329999999999
$ws.Range("H4").Value = 2653.267
$ws.Range("C5").Value = 119.7103015934359
$ws.Range("H5").Value = 2317.934366628473
$ws.Range("C6").Value = 95.69277464312188
$ws.Range("H6").Value = 2053.917872188321
$ws.Range("C7").Value = 71.68043207769836
$ws.Range("H7").Value = 1841.244964503114
$ws.Range("C8").Value = 49.75742254512136
$ws.Range("H8").Value = 1692.806943814478
$ws.Range("C9").Value = 63.88522952843584
$ws.Range("H9").Value = 1738.8476760207
$ws.Range("C10").Value = 80.17042956814157
$ws.Range("H10").Value = 1756.293963726893
$ws.Range("C11").Value = 96.99510127094811
$ws.Range("H11").Value = 1755.554020091888
$ws.Range("C12").Value = 112.932368446504
$ws.Range("H12").Value = 1729.631979650156
$ws.Range("C13").Value = 127.8342527553297
$ws.Range("H13").Value = 1731.529779708525
$ws.Range("C14").Value = 142.4876637965252
$ws.Range("H14").Value = 1734.478848120219
$ws.Range("C15").Value = 157.9219908217015
$ws.Range("H15").Value = 1738.146878041396
$ws.Range("C16").Value = 173.2853254846975
$ws.Range("H16").Value = 1743.874878588342
$ws.Range("C17").Value = 172.7758566111596
$ws.Range("H17").Value = 1761.524043094089
$ws.Range("C18").Value = 172.6923419107023
$ws.Range("H18").Value = 1778.481616477271
$ws.Range("C19").Value = 172.999285202235
$ws.Range("H19").Value = 1792.214568108853
$ws.Range("C20").Value = 173.8741673912084
$ws.Range("H20").Value = 1803.970925515257
$ws.Range("C21").Value = 172.6742687666344
$ws.Range("H21").Value = 1845.574035042185
$ws.Range("C22").Value = 171.4033777798804
$ws.Range("H22").Value = 1880.670991984249
$ws.Range("C23").Value = 170.4164562418464
$ws.Range("H23").Value = 1910.004375292176
$ws.Range("C24").Value = 169.7490003336229
$ws.Range("H24").Value = 1932.953941083782
$ws.Range("C25").Value = 168.9450237981079
$ws.Range("H25").Value = 1945.297494348133
$ws.Range("C26").Value = 168.0345587193225
$ws.Range("H26").Value = 1959.254032082428
$ws.Range("C27").Value = 167.6565363568878
$ws.Range("H27").Value = 1976.123083872504
$ws.Range("C28").Value = 167.4559948999033
$ws.Range("H28").Value = 1995.256001116345
$ws.Range("C29").Value = 166.9802474902943
$ws.Range("H29").Value = 2011.893696146629
$ws.Range("C30").Value = 166.3980115374151
$ws.Range("H30").Value = 2028.836695249428
$ws.Range("C31").Value = 165.993256489986
$ws.Range("H31").Value = 2044.049055535957
$ws.Range("C32").Value = 165.553005261467
$ws.Range("H32").Value = 2058.808439852786
$ws.Range("C33").Value = 164.7961863009631
$ws.Range("H33").Value = 2098.553659222418
$ws.Range("C34").Value = 164.1813520648193
$ws.Range("H34").Value = 2133.44714496436
$ws.Range("C35").Value = 163.5665178286756
$ws.Range("H35").Value = 2164.35538491337
$ws.Range("C36").Value = 162.9161874114417
$ws.Range("H36").Value = 2190.577573660407
$ws.Range("C37").Value = 161.6183775766845
$ws.Range("H37").Value = 2163.361106057697
$ws.Range("C38").Value = 160.320567741927
$ws.Range("H38").Value = 2147.650682851897
$ws.Range("C39").Value = 159.0937502693497
$ws.Range("H39").Value = 2144.846927132937
$ws.Range("C40").Value = 157.724948072412
$ws.Range("H40").Value = 2154.252842106352
$ws.Range("C41").Value = 155.2284846079838
$ws.Range("H41").Value = 2221.525875455438
$ws.Range("C42").Value = 152.8740058679153
$ws.Range("H42").Value = 2278.407505250257
$ws.Range("C43").Value = 150.3420462223971
$ws.Range("H43").Value = 2319.479494032694
$ws.Range("C44").Value = 147.8810789390589
$ws.Range("H44").Value = 2347.944425713181
